$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.168.93"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "1.854.65"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +1.10%  "

$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "310.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3705"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "

$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07818"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").Value = "1.815.64"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.396"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.020"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008714"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("D20").Value = "27.209.73"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.082"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.943"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.307"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.559"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7374"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.684"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.03%  "

$ws.Range("E36").Value = "  +2.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05246"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.058"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1530"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.328"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.20%  "

$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4756"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.017"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.627"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06069"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "

